# Natmi following Dr Hou advice
# Cd44-Sele LR-pairs sheet: recompute rows for all (Sending cluster) x (Target cluster in {ECs, M2}) pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Cd44", "Sele", "ECs", 3, 1, 21.084959, 63.25487699999999, 0.03381320693734752, 0.03509122472428063, 3, 1, 9.112632333333332, 27.337897, 0.9981738658344552, 0.9981738658344552, 192.1394791304076, 1729.255312173669, 0.03375145948491259, 0.03502714343990081),
    @("ECs", "Cd44", "Sele", "M2", 3, 1, 21.084959, 63.25487699999999, 0.03381320693734752, 0.03509122472428063, 1, 0.3333333333333333, 0.01667133333333333, 0.050014, 0.001826134165544791, 0.001826134165544791, 0.3515143798086666, 3.163629418278, 0.00006174745243492645, 0.00006408128437981894),
    @("FAPs", "Cd44", "Sele", "ECs", 3, 1, 81.06331633333333, 243.189949, 0.12999838843446, 0.1349118606466557, 3, 1, 9.112632333333332, 27.337897, 0.9981738658344552, 0.9981738658344552, 738.7001974663614, 6648.301777197253, 0.1297609939358741, 0.1346654934885916),
    @("FAPs", "Cd44", "Sele", "M2", 3, 1, 81.06331633333333, 243.189949, 0.12999838843446, 0.1349118606466557, 1, 0.3333333333333333, 0.01667133333333333, 0.050014, 0.001826134165544791, 0.001826134165544791, 1.351433567698444, 12.162902109286, 0.0002373944985859302, 0.0002463671580640757),
    @("M1", "Cd44", "Sele", "ECs", 3, 1, 136.9994176666667, 410.998253, 0.2197011461990087, 0.2280050605000741, 3, 1, 9.112632333333332, 27.337897, 0.9981738658344552, 0.9981738658344552, 1248.425323077104, 11235.82790769394, 0.2192999424297253, 0.2275886926691778),
    @("M1", "Cd44", "Sele", "M2", 3, 1, 136.9994176666667, 410.998253, 0.2197011461990087, 0.2280050605000741, 1, 0.3333333333333333, 0.01667133333333333, 0.050014, 0.001826134165544791, 0.001826134165544791, 2.283962958393555, 20.555666625542, 0.0004012037692833609, 0.0004163678308962924),
    @("M2", "Cd44", "Sele", "ECs", 3, 1, 316.292811, 948.878433, 0.5072276531881493, 0.5263990368430604, 3, 1, 9.112632333333332, 27.337897, 0.9981738658344552, 0.9981738658344552, 2882.260096319488, 25940.3408668754, 0.5063013874409533, 0.5254377615771715),
    @("M2", "Cd44", "Sele", "M2", 3, 1, 316.292811, 948.878433, 0.5072276531881493, 0.5263990368430604, 1, 0.3333333333333333, 0.01667133333333333, 0.050014, 0.001826134165544791, 0.001826134165544791, 5.273022883117999, 47.457205948062, 0.0009262657471959839, 0.0009612752658889839),
    @("sCs", "Cd44", "Sele", "ECs", 2, 1, 68.131198, 136.262396, 0.1092596052410345, 0.07559281728592908, 3, 1, 9.112632333333332, 27.337897, 0.9981738658344552, 0.9981738658344552, 620.8545578035353, 3725.127346821212, 0.1090600825429899, 0.07545477465961346),
    @("sCs", "Cd44", "Sele", "M2", 2, 1, 68.131198, 136.262396, 0.1092596052410345, 0.07559281728592908, 1, 0.3333333333333333, 0.01667133333333333, 0.050014, 0.001826134165544791, 0.001826134165544791, 1.135837912257333, 6.815027473544, 0.0001995226980445898, 0.0001380426263156199)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
